$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) EDUCATION / UNLV line: "Expected: May of 2024" -> "Expected: Fall of 2024"
#    The target markup splits this into four discrete runs (extra spacing,
#    "Expected: ", "Fall", " of 2024 ") in place of the old single run, and
#    drops one of the two tab runs that used to precede it.
# ---------------------------------------------------------------------------
$rng = $d.Content.Duplicate
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("Expected: May of 2024 ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Pull in the tab character immediately before "Expected:" too, so it
    # can be swapped out for the run of plain spaces the target uses.
    $target = $d.Range($rng.Start - 1, $rng.End)
    $target.Text = ""

    $insertPoint = $d.Range($target.Start, $target.Start)
    $insertPoint.InsertAfter("               ")

    $insertPoint = $d.Range($insertPoint.End, $insertPoint.End)
    $insertPoint.InsertAfter("Expected: ")

    $insertPoint = $d.Range($insertPoint.End, $insertPoint.End)
    $insertPoint.InsertAfter("Fall")

    $insertPoint = $d.Range($insertPoint.End, $insertPoint.End)
    $insertPoint.InsertAfter(" of 2024 ")
}

# ---------------------------------------------------------------------------
# 2) Re-saving the document also materialises the built-in "FollowedHyperlink"
#    character style in styles.xml (Word always pairs it with "Hyperlink").
#    Touch it via an existing hyperlink run and revert the run immediately,
#    which mints the style definition without altering visible content.
# ---------------------------------------------------------------------------
$hrng = $d.Content.Duplicate
$hrng.Find.Execute("krisacuna.com", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($hrng.Find.Found) {
    $origStyle = $hrng.Style
    $hrng.Style = "FollowedHyperlink"
    $hrng.Style = $origStyle

    $fh = $d.Styles("FollowedHyperlink")
    $fh.Priority = 99
    $fh.UnhideWhenUsed = $true
    $fh.QuickStyle = $false
}
